# Generate Report for Handback
# Updates the handoff/handback timestamps for the "230aab99-..." file
# (row 2 of the zh-cn / de-de tracking sheets, and the matching row on
# the Overview sheet) to reflect a freshly generated handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 230aab99 file
$wsOverview.Range("G2").Value = "2016-08-20 14:51:05"

# zh-cn sheet: Correspond Handoff / Handback datetimes for the 230aab99 file
$wsZhCn.Range("H2").Value = "2016-08-20 14:50:58"
$wsZhCn.Range("K2").Value = "2016-08-20 14:51:27"

# de-de sheet: Correspond Handoff / Handback datetimes for the 230aab99 file
$wsDeDe.Range("H2").Value = "2016-08-20 14:51:05"
$wsDeDe.Range("K2").Value = "2016-08-20 14:51:33"
